$d = $word.ActiveDocument

# --- Step 1: Turn the "Meta description" paragraph (currently paragraph 2,
# right under the H1 title) into the new bold "Play Desperados..." paragraph
# that belongs at the bottom of the document (just above the "Prompt:" blurb).
$metaPara = $d.Paragraphs.Item(2)

# Change the bold "Meta description" run into the new bold heading text.
$metaPara.Range.Find.Execute("Meta description", $false, $false, $false, `
    $false, $false, $true, 1, $false, `
    "Play Desperados Wild Megaways slot game for free", 2)

# Remove the trailing plain-text run (": Read our neutral review ...").
$metaPara.Range.Find.Execute(`
    ": Read our neutral review of Desperados Wild Megaways slot game. Play it for free without downloading.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 2: Move that (now 2-run, bold) paragraph down to sit right before
# the final "Prompt: ..." paragraph.
$metaPara.Range.Cut()

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPoint = $lastPara.Range.Start
$target = $d.Range($insertPoint, $insertPoint)
$target.Paste()

# --- Step 3: Replace the italic "Prompt: ..." paragraph's text with the new
# "Read our neutral review ..." text (keeps its existing italic formatting).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute(`
    "Prompt: Create a feature image for Desperados Wild Megaways that is fitting for the Wild West theme of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be standing in front of a saloon, with his arms up in victory and a pile of gold coins behind him. In the background, the Grand Canyon and a Wild West town should be visible. The image should convey the excitement and adventure of the game, while also highlighting the unique character of the Maya warrior.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Read our neutral review of Desperados Wild Megaways slot game. Play it for free without downloading.", 2)
